# Replace the two-digit multiplication problems in the document's table
# with the new set of problems, per the commit's regenerated output.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "83×39="; new = "62×59=" },
    @{ old = "13×79="; new = "84×40=" },
    @{ old = "47×95="; new = "90×89=" },
    @{ old = "83×56="; new = "54×74=" },
    @{ old = "65×15="; new = "33×31=" },
    @{ old = "76×23="; new = "25×98=" },
    @{ old = "22×62="; new = "75×11=" },
    @{ old = "61×86="; new = "81×62=" },
    @{ old = "92×33="; new = "58×70=" },
    @{ old = "23×47="; new = "11×23=" },
    @{ old = "16×61="; new = "42×71=" },
    @{ old = "91×87="; new = "49×27=" },
    @{ old = "15×31="; new = "29×79=" },
    @{ old = "82×42="; new = "24×16=" },
    @{ old = "93×12="; new = "87×94=" },
    @{ old = "23×34="; new = "56×82=" },
    @{ old = "31×32="; new = "54×70=" },
    @{ old = "73×73="; new = "85×99=" },
    @{ old = "44×70="; new = "97×75=" },
    @{ old = "38×74="; new = "76×13=" },
    @{ old = "35×41="; new = "17×24=" },
    @{ old = "19×53="; new = "14×94=" },
    @{ old = "86×51="; new = "60×45=" },
    @{ old = "43×63="; new = "91×74=" },
    @{ old = "76×54="; new = "37×34=" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
